$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "SBO_DEF" in column F1, matching the formatting of the other header cells (e.g. E1)
$ws.Range("F1").Value = "SBO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill F2:F8 with the literal string "[]" (plain data cells, same as columns C/E - no special style)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}

$excel.CutCopyMode = 0
